$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$values = @{
    15 = '2023/2024, 2025/2026'
    19 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    21 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    22 = 'Eman_mohamed@med.asu.edu.eg, 2025/2026'
    37 = '2023/2024, 2025/2026'
    41 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    43 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    44 = 'Eman_mohamed@med.asu.edu.eg, 2025/2026'
    60 = '2026/2027, 2025/2026'
    63 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    64 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    65 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    66 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    82 = '2026/2027, 2025/2026'
    85 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    86 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    87 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    88 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    89 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    104 = '2026/2027, 2025/2026'
    106 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    107 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    110 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    111 = 'neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, 2025/2026'
    126 = '2026/2027, 2025/2026'
    128 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    129 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    132 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    133 = 'neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, 2025/2026'
    150 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    153 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    155 = 'neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, 2025/2026'
    172 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    175 = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
    177 = 'neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, 2025/2026'
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}
